$wb = $excel.ActiveWorkbook

# Update JEMINE sheet's "META AN DISTRI" row (row 6) values from 25000 to 35000
# and move its selection to F10.
$wsJemine = $wb.Worksheets.Item("JEMINE")
$wsJemine.Range("C6:N6").Value = 35000
$wsJemine.Range("F10").Select()

# Update PAOLA sheet's "META AN DISTRI" row (row 6) values from 25000 to 35000
# and move its selection to E23. PAOLA is edited/selected last so it remains
# the active sheet/tab, matching the original workbook's active tab.
$wsPaola = $wb.Worksheets.Item("PAOLA")
$wsPaola.Range("C6:N6").Value = 35000
$wsPaola.Range("E23").Select()
